$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

$ws.Range("N8").Value = "POS800-M"
$ws.Range("P8").Value = "POS800-M - 1"
$ws.Range("Q8").Value = 0.27200000000000002
$ws.Range("R8").Value = 0.27200000000000002

$ws.Range("A8").Select()
